$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: numeric-looking price strings that must stay as literal text
# (matches the workbook's original inlineStr formatting, e.g. trailing zeros / thousand dots)

$ws.Range("D2").Value = "42.756.13"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.546.70"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.26"
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.47"
$ws.Range("E6").Value = "  +6.00%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.42"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.41"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "2.939.32"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.88"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "2.533.95"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "42.792.47"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.41"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0957"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.23"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.91"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.59"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.02"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.27"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.74"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0804"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.63"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("E37").Value = "  +6.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.26"
$ws.Range("E38").Value = "  -5.70%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("E41").Value = "  +10.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.40"
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").Value = "1.982.67"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.00"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "2.793.89"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.44"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("E51").Value = "  -1.53%  "
